$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103, shifting rows 103:172 down to 104:173.
$ws.Rows(103).Insert()

# Fill the new row 103 with the new data.
$ws.Range("A103").Value = 3
$ws.Range("B103").Value = "Femacal de La Calera"
$ws.Range("C103").Value = "Coquimbo"
$ws.Range("D103").Value = 44438
$ws.Range("E103").Value = 5
$ws.Range("F103").Value = 100112012
$ws.Range("G103").Value = "Espinaca"
$ws.Range("H103").Value = "Sin especificar"
$ws.Range("I103").Value = "Primera"
$ws.Range("J103").Value = 120
$ws.Range("K103").Value = 3500
$ws.Range("L103").Value = 3500
$ws.Range("M103").Value = 3500
$ws.Range("N103").Value = '$/docena de atados (3 kilos)'
$ws.Range("O103").Value = "Provincia de Quillota"
$ws.Range("P103").Value = 1167
$ws.Range("Q103").Value = 3
$ws.Range("R103").Value = "Hortaliza"
